# Regenerate the "K" column (column G, header "K", previously "Strike#")
# with freshly calculated values (std/mean-derived s_vals), rows 2-55.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newK = @(0,2,0,0,0,0,1,2,0,1,0,0,1,0,1,0,3,2,0,0,1,1,0,0,0,0,0,2,2,2,2,4,0,1,0,2,1,3,2,2,0,4,2,1,2,0,0,0,2,2,0,1,1,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
